$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "28.258.89"
$ws.Range("E2").Value = "  -2.56%  "

$ws.Range("D3").Value = "1.868.69"
$ws.Range("E3").Value = "  -2.08%  "

Set-TextValue "D4" "1.005"
$ws.Range("E4").Value = "  +0.23%  "

Set-TextValue "D5" "317.92"
$ws.Range("E5").Value = "  -2.04%  "

$ws.Range("E6").Value = "  +0.18%  "

Set-TextValue "D7" "0.4391"
$ws.Range("E7").Value = "  -4.38%  "

Set-TextValue "D8" "0.3685"
$ws.Range("E8").Value = "  -3.56%  "

Set-TextValue "D9" "0.07482"
$ws.Range("E9").Value = "  -3.01%  "

Set-TextValue "D10" "0.9348"
$ws.Range("E10").Value = "  -4.58%  "

Set-TextValue "D11" "21.33"
$ws.Range("E11").Value = "  -3.35%  "

$ws.Range("D12").Value = "1.908.15"
$ws.Range("E12").Value = "  +0.80%  "

Set-TextValue "D13" "6.690"
$ws.Range("E13").Value = "  -3.50%  "

$ws.Range("E14").Value = "  -3.73%  "

Set-TextValue "D15" "0.06896"
$ws.Range("E15").Value = "  -1.88%  "

$ws.Range("E16").Value = "  +0.18%  "

Set-TextValue "D17" "81.98"

Set-TextValue "D18" "0.000009016"
$ws.Range("E18").Value = "  -4.73%  "

$ws.Range("E19").Value = "  +0.16%  "

Set-TextValue "D20" "15.90"
$ws.Range("E20").Value = "  -4.78%  "

$ws.Range("D21").Value = "28.237.29"
$ws.Range("E21").Value = "  -2.48%  "

Set-TextValue "D22" "5.110"
$ws.Range("E22").Value = "  -3.88%  "

Set-TextValue "D23" "10.81"
$ws.Range("E23").Value = "  -0.59%  "

$ws.Range("D24").Value = "2.103.42"
$ws.Range("E24").Value = "  -1.26%  "

$ws.Range("E25").Value = "  -3.15%  "

$ws.Range("E26").Value = "  -1.95%  "

Set-TextValue "D27" "18.39"
$ws.Range("E27").Value = "  -3.47%  "

Set-TextValue "D28" "5.318"
$ws.Range("E28").Value = "  -5.97%  "

Set-TextValue "D29" "113.48"
$ws.Range("E29").Value = "  -3.39%  "

$ws.Range("E30").Value = "  -7.10%  "

Set-TextValue "D31" "0.09034"
$ws.Range("E31").Value = "  -2.70%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D32" "0.7923"
$ws.Range("E32").Value = "  -8.22%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D33" "4.839"
$ws.Range("E33").Value = "  -4.51%  "

Set-TextValue "D34" "1.170"
$ws.Range("E34").Value = "  -6.16%  "

Set-TextValue "D35" "2.926"
$ws.Range("E35").Value = "  -3.30%  "

Set-TextValue "D37" "1.125"
$ws.Range("E37").Value = "  -2.54%  "

Set-TextValue "D38" "0.05439"
$ws.Range("E38").Value = "  -5.31%  "

Set-TextValue "D39" "0.01969"
$ws.Range("E39").Value = "  -3.48%  "

Set-TextValue "D40" "2.955"
$ws.Range("E40").Value = "  +3.94%  "

Set-TextValue "D41" "0.5248"
$ws.Range("E41").Value = "  -4.74%  "

Set-TextValue "D42" "7.058"
$ws.Range("E42").Value = "  -4.70%  "

Set-TextValue "D43" "0.1679"
$ws.Range("E43").Value = "  -4.36%  "

Set-TextValue "D44" "8.696"
$ws.Range("E44").Value = "  -6.72%  "

Set-TextValue "D45" "0.06753"

Set-TextValue "D46" "0.4862"
$ws.Range("E46").Value = "  -6.16%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D47" "10.58"
$ws.Range("E47").Value = "  -5.71%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D48" "106.90"
$ws.Range("E48").Value = "  -3.59%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D49" "1.003"
$ws.Range("E49").Value = "  +0.17%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D50" "1.906"
$ws.Range("E50").Value = "  -6.56%  "

Set-TextValue "D51" "1.668"
$ws.Range("E51").Value = "  -6.24%  "
